$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-14: update Price (column D) only ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '269.75'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.89'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '6.366'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06221'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.633'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.715'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.399'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8356'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01368'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1615'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08292'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03422'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03190'

# --- Rows 15-26: coin list shifted by one position (Coin, Link, Price, Volume(1h)) ---
# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09313'
$ws.Range("E15").Value = '14BitMartTokenBMX'
# Row 16
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.934'
$ws.Range("E16").Value = '15MCDexMCB'
# Row 17
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001711'
$ws.Range("E17").Value = '16BitForexTokenBF'
# Row 18
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04871'
$ws.Range("E18").Value = '17CoinExTokenCET'
# Row 19
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006231'
$ws.Range("E19").Value = '18TigerCashTCH'
# Row 20
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.005391'
$ws.Range("E20").Value = '19HotbitTokenHTB'
# Row 21
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001096'
$ws.Range("E21").Value = '20BitKanKAN'
# Row 22
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001509'
$ws.Range("E22").Value = '21NitroExNTX'
# Row 23
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.762'
$ws.Range("E23").Value = '22LEOLEO'
# Row 24
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.368'
$ws.Range("E24").Value = '23BTSETokenBTSE'
# Row 25
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3333'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
# Row 26
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1248'
$ws.Range("E26").Value = '25ProBitTokenPROB'

# --- Rows 40-50: update Price (column D) only ---
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04673'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006943'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1162'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003483'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01230'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006285'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000755'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.7044'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1398'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002113'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01248'
